# Update crypto price/volume figures to reflect the latest scrape
# Values are stored as text (not numbers) to match the source data,
# so a leading apostrophe is used to force text entry, then the style
# is reset to "Normal" so no stray number-format is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'44.687.24"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +3.54%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.424.52"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.22%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.00%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'315.46"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +3.83%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'101.86"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +6.38%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.517"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +2.73%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = "'0.528"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +9.70%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'35.55"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +3.46%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0801"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.84%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'18.90"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.95%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.11%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'6.96"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +3.13%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.801.59"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.36%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.399.96"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.13%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.834"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +4.47%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'44.523.39"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +3.19%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'12.35"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +2.85%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.20%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.71%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'68.76"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.88%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'242.31"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.95%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  +3.97%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +2.25%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D27').Value = "'25.21"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = "'  -3.47%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +1.66%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'33.65"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +4.04%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'48.55"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.59%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +18.36%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'19.56"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +10.92%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.0778"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +8.46%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +3.09%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.26%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.49%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'4.49"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +2.95%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.76%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'122.55"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.57%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +1.89%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -3.12%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'21.12"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.16%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.0291"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +4.48%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'1.949.16"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.55%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.06%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'2.96"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +8.42%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.99%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.71"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +12.48%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'75.32"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +4.72%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'54.29"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +5.73%  "
$ws.Range('E51').Style = 'Normal'
